$wb = $excel.ActiveWorkbook

# --- Add the new "oracle-resource" sheet after the last existing sheet (macos) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "oracle-resource"

$oracleUrl = "https://docs.oracle.com/en/java/javase/17/install/overview-jdk-installation.html"
$cell = $newSheet.Range("A1")
$cell.Value = $oracleUrl
$newSheet.Hyperlinks.Add($cell, $oracleUrl)

# Match the author's final selection on the new sheet
$newSheet.Range("A4").Select()

# Column width for the new sheet (roughly matches the width used for the resource link column)
$newSheet.Columns.Item(1).ColumnWidth = 67.5
